$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.904.85"
$ws.Range("E2").Value = "  +2.55%  "

$ws.Range("D3").Value = "3.089.85"
$ws.Range("E3").Value = "  +5.20%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.91"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.27%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.60"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +6.23%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").Value = "3.085.07"
$ws.Range("E8").Value = "  +5.17%  "

$ws.Range("E9").Value = "  +1.33%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.64"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.39%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.155"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.53%  "

$ws.Range("E12").Value = "  +5.96%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000250"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.96%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.45"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +6.30%  "

$ws.Range("E15").Value = "  -0.48%  "

$ws.Range("D16").Value = "3.599.26"

$ws.Range("D17").Value = "66.841.29"
$ws.Range("E17").Value = "  +2.43%  "

$ws.Range("E18").Value = "  +3.79%  "

$ws.Range("D19").Value = "3.089.04"
$ws.Range("E19").Value = "  +5.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.15"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +8.38%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "466.97"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.96%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.715"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +4.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.53"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +4.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.65"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.05%  "

$ws.Range("E25").Value = "  +6.85%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.10"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +8.37%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.14"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.48%  "

$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.00"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.19%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.40"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.18%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.67"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +4.14%  "

$ws.Range("E32").Value = "  +1.62%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.26"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +4.39%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.115"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.83%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.03%  "

$ws.Range("E36").Value = "  +3.21%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.89"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +3.29%  "

$ws.Range("B38").Value = "Arweave"
$ws.Range("C38").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "47.00"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +5.59%  "

$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.11"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +7.09%  "

$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.318"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +6.90%  "

$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.23"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.26%  "

$ws.Range("E42").Value = "  +2.02%  "

$ws.Range("E43").Value = "  +2.86%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.83"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0361"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.81%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "382.92"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.53%  "

$ws.Range("D47").Value = "2.774.83"
$ws.Range("E47").Value = "  +2.59%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "135.11"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.44%  "

$ws.Range("E49").Value = "  +0.07%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.82"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +6.65%  "

$ws.Range("E51").Value = "  +1.92%  "
